$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46037
$ws.Range("D8").Value = 152.67
$ws.Range("E8").Value = 149.22
$ws.Range("F8").Value = 159.22
$ws.Range("G8").Value = 149.23
$ws.Range("A9").Value = 46037
$ws.Range("D9").Value = 152.67
$ws.Range("E9").Value = 149.22
$ws.Range("F9").Value = 159.22
$ws.Range("G9").Value = 149.23
$ws.Range("A10").Value = 46037
$ws.Range("D10").Value = 154.32
$ws.Range("E10").Value = 151.72
$ws.Range("F10").Value = 161.72
$ws.Range("G10").Value = 152.12
$ws.Range("A11").Value = 46036
$ws.Range("D11").Value = 153.29
$ws.Range("E11").Value = 149.81
$ws.Range("F11").Value = 159.81
$ws.Range("G11").Value = 149.83
$ws.Range("A12").Value = 46036
$ws.Range("D12").Value = 153.29
$ws.Range("E12").Value = 149.81
$ws.Range("F12").Value = 159.81
$ws.Range("G12").Value = 149.83
$ws.Range("A13").Value = 46036
$ws.Range("D13").Value = 154.4
$ws.Range("E13").Value = 151.43
$ws.Range("F13").Value = 161.43
$ws.Range("G13").Value = 151.83
$ws.Range("A17").Value = 46037
$ws.Range("D17").Value = 158.11
$ws.Range("E17").Value = 154.37
$ws.Range("F17").Value = 164.37
$ws.Range("A18").Value = 46036
$ws.Range("D18").Value = 157.66
$ws.Range("E18").Value = 154.11
$ws.Range("F18").Value = 164.11
$ws.Range("A22").Value = 46037
$ws.Range("D22").Value = 154.08
$ws.Range("E22").Value = 151.23
$ws.Range("F22").Value = 160.83
$ws.Range("G22").Value = 152.3
$ws.Range("A23").Value = 46037
$ws.Range("D23").Value = 159.29
$ws.Range("E23").Value = 157.06
$ws.Range("F23").Value = 167.06
$ws.Range("A24").Value = 46037
$ws.Range("D24").Value = 159.45
$ws.Range("E24").Value = 157.7
$ws.Range("F24").Value = 167.7
$ws.Range("A25").Value = 46037
$ws.Range("D25").Value = 159.43
$ws.Range("E25").Value = 157.23
$ws.Range("F25").Value = 167.23
$ws.Range("G25").Value = 157.36
$ws.Range("A26").Value = 46037
$ws.Range("D26").Value = 159.06
$ws.Range("E26").Value = 158.84
$ws.Range("F26").Value = 168.84
$ws.Range("A27").Value = 46036
$ws.Range("D27").Value = 154.26
$ws.Range("E27").Value = 150.94
$ws.Range("F27").Value = 160.54
$ws.Range("G27").Value = 152.02
$ws.Range("A28").Value = 46036
$ws.Range("D28").Value = 158.82
$ws.Range("E28").Value = 156.77
$ws.Range("F28").Value = 166.77
$ws.Range("A29").Value = 46036
$ws.Range("D29").Value = 158.98
$ws.Range("E29").Value = 157.41
$ws.Range("F29").Value = 167.41
$ws.Range("A30").Value = 46036
$ws.Range("D30").Value = 158.96
$ws.Range("E30").Value = 156.93
$ws.Range("F30").Value = 166.93
$ws.Range("G30").Value = 157.06
$ws.Range("A31").Value = 46036
$ws.Range("D31").Value = 158.59
$ws.Range("E31").Value = 158.54
$ws.Range("F31").Value = 168.54
$ws.Range("A35").Value = 46037
$ws.Range("D35").Value = 152.21
$ws.Range("E35").Value = 148.7
$ws.Range("F35").Value = 157.7
$ws.Range("A36").Value = 46036
$ws.Range("D36").Value = 152.84
$ws.Range("E36").Value = 148.42
$ws.Range("F36").Value = 157.42
$ws.Range("A40").Value = 46037
$ws.Range("D40").Value = 158.84
$ws.Range("E40").Value = 157.1
$ws.Range("F40").Value = 167.1
$ws.Range("A41").Value = 46037
$ws.Range("D41").Value = 158.56
$ws.Range("E41").Value = 157.52
$ws.Range("F41").Value = 167.52
$ws.Range("A42").Value = 46036
$ws.Range("D42").Value = 158.34
$ws.Range("E42").Value = 156.74
$ws.Range("F42").Value = 166.74
$ws.Range("A43").Value = 46036
$ws.Range("D43").Value = 158.06
$ws.Range("E43").Value = 157.16
$ws.Range("F43").Value = 167.16
$ws.Range("A47").Value = 46037
$ws.Range("D47").Value = 152.34
$ws.Range("E47").Value = 150.18
$ws.Range("F47").Value = 160.18
$ws.Range("A48").Value = 46037
$ws.Range("D48").Value = 151.97
$ws.Range("E48").Value = 150.11
$ws.Range("F48").Value = 160.11
$ws.Range("A49").Value = 46036
$ws.Range("D49").Value = 152.82
$ws.Range("E49").Value = 150.05
$ws.Range("F49").Value = 160.05
$ws.Range("A50").Value = 46036
$ws.Range("D50").Value = 152.46
$ws.Range("E50").Value = 149.98
$ws.Range("F50").Value = 159.98
$ws.Range("A54").Value = 46037
$ws.Range("D54").Value = 168.03
$ws.Range("E54").Value = 164.32
$ws.Range("F54").Value = 174.32
$ws.Range("A55").Value = 46037
$ws.Range("D55").Value = 160.49
$ws.Range("E55").Value = 162.44
$ws.Range("F55").Value = 172.44
$ws.Range("A56").Value = 46037
$ws.Range("D56").Value = 156.91
$ws.Range("A57").Value = 46037
$ws.Range("D57").Value = 157.51
$ws.Range("E57").Value = 156.86
$ws.Range("A58").Value = 46037
$ws.Range("D58").Value = 153.28
$ws.Range("E58").Value = 152.76
$ws.Range("F58").Value = 162.76
$ws.Range("A59").Value = 46037
$ws.Range("D59").Value = 160.46
$ws.Range("E59").Value = 162.65
$ws.Range("A60").Value = 46036
$ws.Range("D60").Value = 167.56
$ws.Range("E60").Value = 163.92
$ws.Range("F60").Value = 173.92
$ws.Range("A61").Value = 46036
$ws.Range("D61").Value = 160.74
$ws.Range("E61").Value = 162.95
$ws.Range("F61").Value = 172.95
$ws.Range("A62").Value = 46036
$ws.Range("D62").Value = 156.98
$ws.Range("A63").Value = 46036
$ws.Range("D63").Value = 157.62
$ws.Range("E63").Value = 157.37
$ws.Range("A64").Value = 46036
$ws.Range("D64").Value = 153.39
$ws.Range("E64").Value = 153.27
$ws.Range("F64").Value = 163.27
$ws.Range("A65").Value = 46036
$ws.Range("D65").Value = 160.02
$ws.Range("E65").Value = 162.3
